# Fruta / hortaliza, semanal
# Insert two new weekly price rows for "Choclo" (Choclero, Terminal Hortofrutícola
# Agro Chillán) at rows 105-106, pushing the existing rows (old 105..205) down to
# (new 107..207).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows before the current row 105 - everything from old row 105
# onward shifts down by two rows.
$ws.Rows("105:106").Insert()

# New row 105: Choclero / Primera, Región del Maule
$ws.Range("A105").Value = 7
$ws.Range("B105").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C105").Value = "Ñuble"
$ws.Range("D105").Value = 44658
$ws.Range("E105").Value = 16
$ws.Range("F105").Value = 100112024
$ws.Range("G105").Value = "Choclo"
$ws.Range("H105").Value = "Choclero"
$ws.Range("I105").Value = "Primera"
$ws.Range("J105").Value = 12000
$ws.Range("K105").Value = 250
$ws.Range("L105").Value = 300
$ws.Range("M105").Value = 275
$ws.Range("N105").Value = "$/unidad"
$ws.Range("O105").Value = "Región del Maule"
$ws.Range("P105").Value = 275
$ws.Range("Q105").Value = 1
$ws.Range("R105").Value = "Hortaliza"

# New row 106: Choclero / Segunda, Región del Maule
$ws.Range("A106").Value = 7
$ws.Range("B106").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C106").Value = "Ñuble"
$ws.Range("D106").Value = 44658
$ws.Range("E106").Value = 16
$ws.Range("F106").Value = 100112024
$ws.Range("G106").Value = "Choclo"
$ws.Range("H106").Value = "Choclero"
$ws.Range("I106").Value = "Segunda"
$ws.Range("J106").Value = 12000
$ws.Range("K106").Value = 180
$ws.Range("L106").Value = 200
$ws.Range("M106").Value = 190
$ws.Range("N106").Value = "$/unidad"
$ws.Range("O106").Value = "Región del Maule"
$ws.Range("P106").Value = 190
$ws.Range("Q106").Value = 1
$ws.Range("R106").Value = "Hortaliza"
